# Revert "update dasign file"
# This reverts commit 2efe7df9b4f0327842d00b452bfdf5d7c2d7f103.
#
# Changes applied:
#  1) The "Design ---" heading paragraph gains a <w:rFonts w:hint="cs"/> on its
#     paragraph mark run properties.
#  2) The (bookmark-only) paragraph right after it loses that same
#     <w:rFonts w:hint="cs"/> hint, and the following empty centered
#     (sz=24) paragraph is removed while the _GoBack bookmark is relocated
#     to sit immediately before the "תיאור המערכת:" run.
#  3) A stray <w:lastRenderedPageBreak/> is dropped from the run that
#     starts the "אחרת " list item later in the document.

$d = $word.ActiveDocument
$W = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# ---------------------------------------------------------------------
# Edit 1: restructure the "Design ---" / bookmark / empty / "teur" block
# ---------------------------------------------------------------------
$startRng = $d.Content
$null = $startRng.Find.Execute("Design --- ---", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startPara = $startRng.Paragraphs(1)

$endRng = $d.Content
$null = $endRng.Find.Execute("תיאור המערכת:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPara = $endRng.Paragraphs(1)

$blockRng = $d.Range($startPara.Range.Start, $endPara.Range.End)

$frag1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + "`r`n" + '<w:document xmlns:w="' + $W + '"><w:body><w:p w:rsidR="00047DE5" w:rsidRDefault="00A73442" w:rsidP="005C64CD"><w:pPr><w:bidi/><w:jc w:val="center"/><w:rPr><w:rFonts w:hint="cs"/><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:rtl/></w:rPr></w:pPr><w:r w:rsidRPr="00047DE5"><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t>D</w:t></w:r><w:r w:rsidR="00047DE5" w:rsidRPr="00047DE5"><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t>esign</w:t></w:r><w:r w:rsidR="005B0587"><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t xml:space="preserve"> ---</w:t></w:r><w:r w:rsidR="005B0587"><w:rPr><w:rFonts w:hint="cs"/><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00047DE5"><w:rPr><w:rFonts w:hint="cs"/><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:rtl/></w:rPr><w:t>--</w:t></w:r><w:r w:rsidR="005B0587"><w:rPr><w:rFonts w:hint="cs"/><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:rtl/></w:rPr><w:t>-</w:t></w:r></w:p><w:p w:rsidR="00EB6922" w:rsidRDefault="00EB6922" w:rsidP="00EB6922"><w:pPr><w:bidi/><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:rtl/></w:rPr></w:pPr></w:p><w:p w:rsidR="00E728E2" w:rsidRDefault="00E728E2" w:rsidP="00533540"><w:pPr><w:bidi/><w:rPr><w:b/><w:bCs/><w:rtl/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:hint="cs"/><w:b/><w:bCs/><w:rtl/></w:rPr><w:t>תיאור המערכת:</w:t></w:r></w:p></w:body></w:document>'
$blockRng.InsertXML($frag1)

# ---------------------------------------------------------------------
# Edit 2: drop <w:lastRenderedPageBreak/> from the "אחרת " list paragraph
# ---------------------------------------------------------------------
$breakRng = $d.Content
$null = $breakRng.Find.Execute("אחרת נעבור על כל הרשומות בטבלה וכל רשומה שהערכים", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$breakPara = $breakRng.Paragraphs(1)

$frag2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + "`r`n" + '<w:document xmlns:w="' + $W + '"><w:body><w:p w:rsidR="00690A05" w:rsidRPr="00820457" w:rsidRDefault="008D49FC" w:rsidP="008D49FC"><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr><w:bidi/><w:spacing w:line="256" w:lineRule="auto"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t xml:space="preserve">אחרת </w:t></w:r><w:r w:rsidR="00690A05"><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t>נעבור על כל הרשומות בטבלה וכל רשומה שהערכים</w:t></w:r><w:r w:rsidR="008B6069"><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00690A05"><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t>(ה</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t>עמודות) שלה עומדים בתנאים נחזיר</w:t></w:r><w:r w:rsidR="00690A05"><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t xml:space="preserve"> אותה</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document>'
$breakPara.Range.InsertXML($frag2)
